$d = $word.ActiveDocument

# --- Step 1: replace the old sentence with the new sentence text (single run) ---
$openQuote = [char]0x201C
$closeQuote = [char]0x201D
$old = "3.A.1 Se llama al CU " + $openQuote + "Buscar paciente" + $closeQuote + " y el mismo confirma que el paciente no existe. "
$new = "3.A.1 El sistema busca el paciente y el mismo no existe. "

$found = $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) {
    throw "Could not find the original 3.A.1 sentence to replace."
}

# --- Step 2: locate the freshly-written sentence again, to get its Start offset ---
$rng = $d.Content
$found2 = $rng.Find.Execute($new, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not re-locate the replaced 3.A.1 sentence."
}
$base = $rng.Start

# The run currently reads as a single run:
#   "3.A.1 El sistema busca el paciente y el mismo no existe. "
# We need it split into three runs with identical rPr:
#   [0,6)   "3.A.1 "
#   [6,34)  "El sistema busca el paciente"
#   [34,57) " y el mismo no existe. "
# Toggling Bold on/off over a sub-range forces Word to break the run at that
# boundary without leaving any lasting formatting difference behind.

$split1 = $base + 6
$split2 = $base + 6 + 28

# Split point 1: [base, split1) vs rest
$r1 = $d.Range($base, $split1)
$r1.Bold = 1
$r1b = $d.Range($base, $split1)
$r1b.Bold = 0

# Split point 2: [split1, split2) vs rest
$r2 = $d.Range($split1, $split2)
$r2.Bold = 1
$r2b = $d.Range($split1, $split2)
$r2b.Bold = 0

Write-Output "done"
